{"js": "// Update the \"NTT Data\" address block on the first page:\n//   \"Cluj Napoca\" (own paragraph) + \"City, State \" & \"400158\" (next\n//   paragraph) become a single paragraph reading\n//   \"Cluj Napoca City, 400158\", and the now-empty \"City, State 400158\"\n//   paragraph is removed entirely.\n\nconst body = context.document.body;\n\n// 1) Merge the city/zip text into the \"Cluj Napoca\" run.\nconst cluj = body.search(\"Cluj Napoca\", { matchCase: true });\ncluj.load(\"text\");\nawait context.sync();\ncluj.items[0].insertText(\"Cluj Napoca City, 400158\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) Remove the old \"City, State 400158\" paragraph completely.\nbody.paragraphs.load(\"text\");\nawait context.sync();\nfor (const p of body.paragraphs.items) {\n  if (p.text === \"City, State 400158\") {\n    p.delete();\n    break;\n  }\n}\nawait context.sync();\n", "ps1": "# Update the \"NTT Data\" address block on the first page:\n#   \"Cluj Napoca\" (own paragraph) + \"City, State \" & \"400158\" (next paragraph)\n# become a single paragraph reading \"Cluj Napoca City, 400158\", and the\n# now-empty \"City, State 400158\" paragraph is removed entirely.\n\n$d = $word.ActiveDocument\n\n# 1) Merge the city/zip text into the \"Cluj Napoca\" run.\n$find = $d.Content.Find\n$find.Text = \"Cluj Napoca\"\n$find.Replacement.Text = \"Cluj Napoca City, 400158\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# 2) Remove the old \"City, State 400158\" paragraph completely.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -eq \"City, State 400158`r\") {\n        $p.Range.Delete()\n        break\n    }\n}\n"}
